$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 274, shifting the existing row 274..341 down to 275..342
# (mirrors the canonical diff which effectively inserts one new weekly record).
$ws.Rows.Item(274).Insert()

# Populate the newly inserted row 274 with the new record's data.
$ws.Range("A274").Value = 5
$ws.Range("B274").Value = "Macroferia Regional de Talca"
$ws.Range("C274").Value = "Maule"
$ws.Range("D274").Value = 44889
$ws.Range("E274").Value = 7
$ws.Range("F274").Value = 100112009
$ws.Range("G274").Value = "Acelga"
$ws.Range("H274").Value = "Sin especificar"
$ws.Range("I274").Value = "Primera"
$ws.Range("J274").Value = 500
$ws.Range("K274").Value = 2000
$ws.Range("L274").Value = 2000
$ws.Range("M274").Value = 2000
$ws.Range("N274").Value = "$/docena de atados (4 kilos)"
$ws.Range("O274").Value = "Región del Maule"
$ws.Range("P274").Value = 500
$ws.Range("Q274").Value = 4
$ws.Range("R274").Value = "Hortaliza"
